# Natmi following Dr Hou advice
# Update Slit2-Robo4 sending/target-cluster matrix with ECs/FAPs/sCs cell types

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit2"
$ws.Range("C2").Value = "Robo4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.143896
$ws.Range("H2").Value = 0.431688
$ws.Range("I2").Value = 0.02807111181859822
$ws.Range("J2").Value = 0.02807111181859822
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 33.790225
$ws.Range("N2").Value = 101.370675
$ws.Range("O2").Value = 0.9910539426277148
$ws.Range("P2").Value = 0.991053942627715
$ws.Range("Q2").Value = 4.8622782166
$ws.Range("R2").Value = 43.7605039494
$ws.Range("S2").Value = 0.02781998604176521
$ws.Range("T2").Value = 0.02781998604176521

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit2"
$ws.Range("C3").Value = "Robo4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.143896
$ws.Range("H3").Value = 0.431688
$ws.Range("I3").Value = 0.02807111181859822
$ws.Range("J3").Value = 0.02807111181859822
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.039185
$ws.Range("N3").Value = 0.117555
$ws.Range("O3").Value = 0.001149280560927517
$ws.Range("P3").Value = 0.001149280560927517
$ws.Range("Q3").Value = 0.005638564760000001
$ws.Range("R3").Value = 0.05074708284
$ws.Range("S3").Value = 0.00003226158313673762
$ws.Range("T3").Value = 0.00003226158313673762

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit2"
$ws.Range("C4").Value = "Robo4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.143896
$ws.Range("H4").Value = 0.431688
$ws.Range("I4").Value = 0.02807111181859822
$ws.Range("J4").Value = 0.02807111181859822
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.265833
$ws.Range("N4").Value = 0.7974990000000001
$ws.Range("O4").Value = 0.007796776811357526
$ws.Range("P4").Value = 0.007796776811357526
$ws.Range("Q4").Value = 0.03825230536800001
$ws.Range("R4").Value = 0.344270748312
$ws.Range("S4").Value = 0.0002188641936962708
$ws.Range("T4").Value = 0.0002188641936962708

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit2"
$ws.Range("C5").Value = "Robo4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.277274333333334
$ws.Range("H5").Value = 12.831823
$ws.Range("I5").Value = 0.8344071140950421
$ws.Range("J5").Value = 0.8344071140950421
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 33.790225
$ws.Range("N5").Value = 101.370675
$ws.Range("O5").Value = 0.9910539426277148
$ws.Range("P5").Value = 0.991053942627715
$ws.Range("Q5").Value = 144.5300621100583
$ws.Range("R5").Value = 1300.770558990525
$ws.Range("S5").Value = 0.826942460180505
$ws.Range("T5").Value = 0.8269424601805051

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit2"
$ws.Range("C6").Value = "Robo4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.277274333333334
$ws.Range("H6").Value = 12.831823
$ws.Range("I6").Value = 0.8344071140950421
$ws.Range("J6").Value = 0.8344071140950421
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.039185
$ws.Range("N6").Value = 0.117555
$ws.Range("O6").Value = 0.001149280560927517
$ws.Range("P6").Value = 0.001149280560927517
$ws.Range("Q6").Value = 0.1676049947516667
$ws.Range("R6").Value = 1.508444952765
$ws.Range("S6").Value = 0.0009589678761290607
$ws.Range("T6").Value = 0.0009589678761290607

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit2"
$ws.Range("C7").Value = "Robo4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.277274333333334
$ws.Range("H7").Value = 12.831823
$ws.Range("I7").Value = 0.8344071140950421
$ws.Range("J7").Value = 0.8344071140950421
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.265833
$ws.Range("N7").Value = 0.7974990000000001
$ws.Range("O7").Value = 0.007796776811357526
$ws.Range("P7").Value = 0.007796776811357526
$ws.Range("Q7").Value = 1.137040667853
$ws.Range("R7").Value = 10.233366010677
$ws.Range("S7").Value = 0.006505686038407978
$ws.Range("T7").Value = 0.006505686038407978

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit2"
$ws.Range("C8").Value = "Robo4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7049536666666666
$ws.Range("H8").Value = 2.114861
$ws.Range("I8").Value = 0.1375217740863597
$ws.Range("J8").Value = 0.1375217740863597
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 33.790225
$ws.Range("N8").Value = 101.370675
$ws.Range("O8").Value = 0.9910539426277148
$ws.Range("P8").Value = 0.991053942627715
$ws.Range("Q8").Value = 23.82054301124166
$ws.Range("R8").Value = 214.384887101175
$ws.Range("S8").Value = 0.1362914964054447
$ws.Range("T8").Value = 0.1362914964054447

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit2"
$ws.Range("C9").Value = "Robo4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7049536666666666
$ws.Range("H9").Value = 2.114861
$ws.Range("I9").Value = 0.1375217740863597
$ws.Range("J9").Value = 0.1375217740863597
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.039185
$ws.Range("N9").Value = 0.117555
$ws.Range("O9").Value = 0.001149280560927517
$ws.Range("P9").Value = 0.001149280560927517
$ws.Range("Q9").Value = 0.02762360942833333
$ws.Range("R9").Value = 0.248612484855
$ws.Range("S9").Value = 0.0001580511016617188
$ws.Range("T9").Value = 0.0001580511016617188

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit2"
$ws.Range("C10").Value = "Robo4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7049536666666666
$ws.Range("H10").Value = 2.114861
$ws.Range("I10").Value = 0.1375217740863597
$ws.Range("J10").Value = 0.1375217740863597
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.265833
$ws.Range("N10").Value = 0.7974990000000001
$ws.Range("O10").Value = 0.007796776811357526
$ws.Range("P10").Value = 0.007796776811357526
$ws.Range("Q10").Value = 0.187399948071
$ws.Range("R10").Value = 1.686599532639
$ws.Range("S10").Value = 0.001072226579253278
$ws.Range("T10").Value = 0.001072226579253278
